{"js": "// Replace each three-digit-division-by-one-digit equation in the table\n// with its new value, matching the author's diff exactly. Every text run\n// in this document is a unique equation string, so an exact (case- and\n// whole-match) search-and-replace per pair is unambiguous.\nconst replacements = [\n  [\"841\u00f74=210, 1\", \"236\u00f75=47, 1\"],\n  [\"289\u00f77=41, 2\", \"537\u00f74=134, 1\"],\n  [\"623\u00f76=103, 5\", \"413\u00f77=59, 0\"],\n  [\"743\u00f79=82, 5\", \"238\u00f78=29, 6\"],\n  [\"873\u00f78=109, 1\", \"968\u00f72=484, 0\"],\n  [\"500\u00f77=71, 3\", \"621\u00f78=77, 5\"],\n  [\"231\u00f78=28, 7\", \"154\u00f74=38, 2\"],\n  [\"402\u00f78=50, 2\", \"376\u00f78=47, 0\"],\n  [\"364\u00f73=121, 1\", \"401\u00f77=57, 2\"],\n  [\"801\u00f76=133, 3\", \"615\u00f73=205, 0\"],\n  [\"458\u00f74=114, 2\", \"347\u00f72=173, 1\"],\n  [\"740\u00f74=185, 0\", \"172\u00f73=57, 1\"],\n  [\"803\u00f76=133, 5\", \"343\u00f79=38, 1\"],\n  [\"471\u00f74=117, 3\", \"695\u00f79=77, 2\"],\n  [\"425\u00f72=212, 1\", \"156\u00f75=31, 1\"],\n  [\"715\u00f79=79, 4\", \"888\u00f76=148, 0\"],\n  [\"424\u00f77=60, 4\", \"182\u00f78=22, 6\"],\n  [\"704\u00f78=88, 0\", \"673\u00f77=96, 1\"],\n  [\"686\u00f74=171, 2\", \"106\u00f78=13, 2\"],\n  [\"162\u00f78=20, 2\", \"975\u00f79=108, 3\"],\n  [\"133\u00f72=66, 1\", \"406\u00f73=135, 1\"],\n  [\"821\u00f73=273, 2\", \"420\u00f76=70, 0\"],\n  [\"963\u00f78=120, 3\", \"408\u00f76=68, 0\"],\n  [\"229\u00f73=76, 1\", \"114\u00f72=57, 0\"],\n  [\"196\u00f78=24, 4\", \"193\u00f76=32, 1\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = context.document.body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (const item of results.items) {\n    item.insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Replace each three-digit-division-by-one-digit equation in the table\n# with its new value, matching the author's diff exactly. Every text run\n# in this document is a unique equation string, so an exact (case-\n# sensitive, non-wildcard) Find/Replace per pair is unambiguous.\n$d = $word.ActiveDocument\n\n$pairs = @(\n  @(\"841\u00f74=210, 1\", \"236\u00f75=47, 1\"),\n  @(\"289\u00f77=41, 2\", \"537\u00f74=134, 1\"),\n  @(\"623\u00f76=103, 5\", \"413\u00f77=59, 0\"),\n  @(\"743\u00f79=82, 5\", \"238\u00f78=29, 6\"),\n  @(\"873\u00f78=109, 1\", \"968\u00f72=484, 0\"),\n  @(\"500\u00f77=71, 3\", \"621\u00f78=77, 5\"),\n  @(\"231\u00f78=28, 7\", \"154\u00f74=38, 2\"),\n  @(\"402\u00f78=50, 2\", \"376\u00f78=47, 0\"),\n  @(\"364\u00f73=121, 1\", \"401\u00f77=57, 2\"),\n  @(\"801\u00f76=133, 3\", \"615\u00f73=205, 0\"),\n  @(\"458\u00f74=114, 2\", \"347\u00f72=173, 1\"),\n  @(\"740\u00f74=185, 0\", \"172\u00f73=57, 1\"),\n  @(\"803\u00f76=133, 5\", \"343\u00f79=38, 1\"),\n  @(\"471\u00f74=117, 3\", \"695\u00f79=77, 2\"),\n  @(\"425\u00f72=212, 1\", \"156\u00f75=31, 1\"),\n  @(\"715\u00f79=79, 4\", \"888\u00f76=148, 0\"),\n  @(\"424\u00f77=60, 4\", \"182\u00f78=22, 6\"),\n  @(\"704\u00f78=88, 0\", \"673\u00f77=96, 1\"),\n  @(\"686\u00f74=171, 2\", \"106\u00f78=13, 2\"),\n  @(\"162\u00f78=20, 2\", \"975\u00f79=108, 3\"),\n  @(\"133\u00f72=66, 1\", \"406\u00f73=135, 1\"),\n  @(\"821\u00f73=273, 2\", \"420\u00f76=70, 0\"),\n  @(\"963\u00f78=120, 3\", \"408\u00f76=68, 0\"),\n  @(\"229\u00f73=76, 1\", \"114\u00f72=57, 0\"),\n  @(\"196\u00f78=24, 4\", \"193\u00f76=32, 1\")\n)\n\nforeach ($pair in $pairs) {\n  $old = $pair[0]\n  $new = $pair[1]\n  $find = $d.Content.Find\n  $find.Text = $old\n  $find.Replacement.Text = $new\n  $find.Execute($old, $false, $false, $false, $false, $false, $true, 1, $false, $new, 2)\n}\n"}
